$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (prices / volume %) scraped on
# Thu Jul 27 03:51:08 UTC 2023. Leading "'" forces text storage so
# numeric-looking values (e.g. "1.000", "0.9999") are not coerced to
# actual numbers by Excel - the source data are text cells.

$ws.Range("D2").Value = "'29.431.90"
$ws.Range("E2").Value = "'  +0.70%  "

$ws.Range("D3").Value = "'1.877.30"
$ws.Range("E3").Value = "'  +1.02%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  +0.06%  "

$ws.Range("D5").Value = "'0.7204"
$ws.Range("E5").Value = "'  +1.42%  "

$ws.Range("E6").Value = "'  +0.86%  "

$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "'  +0.08%  "

$ws.Range("D8").Value = "'0.07835"
$ws.Range("E8").Value = "'  -2.21%  "

$ws.Range("D9").Value = "'0.3108"
$ws.Range("E9").Value = "'  +2.43%  "

$ws.Range("D10").Value = "'24.98"
$ws.Range("E10").Value = "'  +6.12%  "

$ws.Range("D11").Value = "'0.08255"
$ws.Range("E11").Value = "'  +0.72%  "

$ws.Range("B12").Value = "'WrappedEther"
$ws.Range("C12").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.872.36"
$ws.Range("E12").Value = "'  +0.89%  "

$ws.Range("B13").Value = "'Polygon"
$ws.Range("C13").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.7276"
$ws.Range("E13").Value = "'  +3.23%  "

$ws.Range("B14").Value = "'Polkadot"
$ws.Range("C14").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.282"
$ws.Range("E14").Value = "'  +2.06%  "

$ws.Range("D15").Value = "'91.36"
$ws.Range("E15").Value = "'  +1.81%  "

$ws.Range("D16").Value = "'29.335.25"
$ws.Range("E16").Value = "'  +0.50%  "

$ws.Range("D17").Value = "'5.930"
$ws.Range("E17").Value = "'  +1.67%  "

$ws.Range("D18").Value = "'245.00"
$ws.Range("E18").Value = "'  +2.86%  "

$ws.Range("D19").Value = "'0.000007890"
$ws.Range("E19").Value = "'  +0.17%  "

$ws.Range("D20").Value = "'13.31"
$ws.Range("E20").Value = "'  +0.26%  "

$ws.Range("D21").Value = "'2.115.39"
$ws.Range("E21").Value = "'  +1.55%  "

$ws.Range("D22").Value = "'0.9988"
$ws.Range("E22").Value = "'  +0.11%  "

$ws.Range("D23").Value = "'7.951"
$ws.Range("E23").Value = "'  +6.83%  "

$ws.Range("D24").Value = "'0.9998"
$ws.Range("E24").Value = "'  +0.05%  "

$ws.Range("D25").Value = "'0.1587"
$ws.Range("E25").Value = "'  +9.69%  "

$ws.Range("D26").Value = "'163.91"
$ws.Range("E26").Value = "'  +0.79%  "

$ws.Range("D27").Value = "'9.037"
$ws.Range("E27").Value = "'  +1.13%  "

$ws.Range("D28").Value = "'18.33"
$ws.Range("E28").Value = "'  +1.32%  "

$ws.Range("D29").Value = "'1.364"
$ws.Range("E29").Value = "'  -4.70%  "

$ws.Range("E30").Value = "'  +0.27%  "

$ws.Range("D31").Value = "'4.389"
$ws.Range("E31").Value = "'  +0.42%  "

$ws.Range("D32").Value = "'4.150"
$ws.Range("E32").Value = "'  +3.37%  "

$ws.Range("D33").Value = "'0.05281"
$ws.Range("E33").Value = "'  +1.26%  "

$ws.Range("D34").Value = "'1.942"
$ws.Range("E34").Value = "'  +0.52%  "

$ws.Range("E35").Value = "'  +3.45%  "

$ws.Range("D36").Value = "'0.7217"
$ws.Range("E36").Value = "'  +1.54%  "

$ws.Range("D37").Value = "'2.676"
$ws.Range("E37").Value = "'  +0.24%  "

$ws.Range("D38").Value = "'0.01868"
$ws.Range("E38").Value = "'  +0.45%  "

$ws.Range("D39").Value = "'1.240.14"
$ws.Range("E39").Value = "'  +9.84%  "

$ws.Range("D40").Value = "'2.721"
$ws.Range("E40").Value = "'  -0.04%  "

$ws.Range("D41").Value = "'0.9081"
$ws.Range("E41").Value = "'  -2.25%  "

$ws.Range("D42").Value = "'73.19"
$ws.Range("E42").Value = "'  +3.91%  "

$ws.Range("D43").Value = "'6.087"
$ws.Range("E43").Value = "'  +3.99%  "

$ws.Range("D45").Value = "'103.60"
$ws.Range("E45").Value = "'  +0.70%  "

$ws.Range("D46").Value = "'0.5328"
$ws.Range("E46").Value = "'  -0.13%  "

$ws.Range("B47").Value = "'RocketPoolETH"
$ws.Range("C47").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "'2.013.53"
$ws.Range("E47").Value = "'  +1.87%  "

$ws.Range("B48").Value = "'BabyDogeCoin"
$ws.Range("C48").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000121"
$ws.Range("E48").Value = "'  +1.21%  "

$ws.Range("E49").Value = "'  +12.77%  "

$ws.Range("B50").Value = "'RenderToken"
$ws.Range("C50").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.764"
$ws.Range("E50").Value = "'  -0.18%  "

$ws.Range("B51").Value = "'TheSandbox"
$ws.Range("C51").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "'0.4332"
$ws.Range("E51").Value = "'  +1.65%  "
